$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-476) holds the "Förändrad" (changed) date for every record.
# The whole column was bumped by one day: 45202 -> 45203 (2023-10-03 -> 2023-10-04).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45203
